$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map the numeric "temp" codes in column A to letter codes a/b/c/d
$map = @{15 = "a"; 20 = "b"; 25 = "c"; 30 = "d"}

for ($r = 2; $r -le 81; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $key = [int]$cell.Value()
    $cell.Value = $map[$key]
}

# Update the view: scroll so row 65 is at the top-left, then select C82
try {
    $excel.ActiveWindow.ScrollRow = 65
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
    # scrolling not supported by host; ignore
}
$ws.Range("C82").Select()
try {
    $excel.Goto($ws.Range("C82"), $true)
} catch {
    # Goto-with-scroll not supported by host; ignore
}
